$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Daniel Gafford / PF,C / Dallas Mavericks -> Jalen Johnson / PF / Atlanta Hawks
$ws.Range("A7").Value = "Jalen Johnson"
$ws.Range("B7").Value = "PF"
$ws.Range("C7").Value = "Atlanta Hawks"

# Row 14: Keegan Murray / SF,PF / Sacramento Kings -> Daniel Gafford / PF,C / Dallas Mavericks
$ws.Range("A14").Value = "Daniel Gafford"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Dallas Mavericks"

# Row 15: Stephen Curry / PG,SG / Golden State Warriors -> Keegan Murray / SF,PF / Sacramento Kings
$ws.Range("A15").Value = "Keegan Murray"
$ws.Range("B15").Value = "SF,PF"
$ws.Range("C15").Value = "Sacramento Kings"

# Row 16: Jalen Johnson / PF / Atlanta Hawks -> Stephen Curry / PG,SG / Golden State Warriors
$ws.Range("A16").Value = "Stephen Curry"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Golden State Warriors"
